# ----------------------------------------------------------------------
# Adds the "curve1" moment-curvature worksheet, wires it up from the
# Sections sheet, extends several input tables (Coordinates, Supports,
# Line Elements, Nodal Load) with extra rows / revised values.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Coordinates sheet: tighten the X-coordinate spacing (0.1 step
#    instead of 0.15) and extend the table with two more nodes.
# ------------------------------------------------------------------
$wsCoord = $wb.Worksheets.Item("Coordinates")

$xVals = @(0.3,0.4,0.5,0.6,0.7,0.8,0.9,1,1.1,1.2,1.3,1.4,1.5,1.6,1.7,1.8)
for ($i = 0; $i -lt $xVals.Length; $i++) {
  $wsCoord.Cells.Item(5 + $i, 2).Value = $xVals[$i]
}

$wsCoord.Cells.Item(21, 1).Value = 20
$wsCoord.Cells.Item(21, 2).Value = 1.9
$wsCoord.Cells.Item(21, 3).Value = 0
$wsCoord.Cells.Item(21, 4).Value = 0

$wsCoord.Cells.Item(22, 1).Value = 21
$wsCoord.Cells.Item(22, 2).Value = 2
$wsCoord.Cells.Item(22, 3).Value = 0
$wsCoord.Cells.Item(22, 4).Value = 0

$wsCoord.Range("D35").Select()

# ------------------------------------------------------------------
# 2) Supports sheet: flip two Tx flags and extend with two more
#    support rows.
# ------------------------------------------------------------------
$wsSupp = $wb.Worksheets.Item("Supports")

$wsSupp.Range("D12").Value = 1
$wsSupp.Range("D20").Value = 0

$wsSupp.Cells.Item(21, 1).Value = 20
$wsSupp.Cells.Item(21, 2).Value = 20
$wsSupp.Cells.Item(21, 3).Value = 1
$wsSupp.Cells.Item(21, 4).Value = 0
$wsSupp.Cells.Item(21, 5).Value = 1
$wsSupp.Cells.Item(21, 6).Value = 1
$wsSupp.Cells.Item(21, 7).Value = 1
$wsSupp.Cells.Item(21, 8).Value = 0

$wsSupp.Cells.Item(22, 1).Value = 21
$wsSupp.Cells.Item(22, 2).Value = 21
$wsSupp.Cells.Item(22, 3).Value = 1
$wsSupp.Cells.Item(22, 4).Value = 1
$wsSupp.Cells.Item(22, 5).Value = 1
$wsSupp.Cells.Item(22, 6).Value = 1
$wsSupp.Cells.Item(22, 7).Value = 1
$wsSupp.Cells.Item(22, 8).Value = 0

$wsSupp.Range("B12").Select()

# ------------------------------------------------------------------
# 3) Line Elements sheet: extend with two more line elements.
# ------------------------------------------------------------------
$wsLine = $wb.Worksheets.Item("Line Elements")

$wsLine.Cells.Item(20, 1).Value = 19
$wsLine.Cells.Item(20, 2).Value = 4
$wsLine.Cells.Item(20, 3).Value = 19
$wsLine.Cells.Item(20, 4).Value = 20
$wsLine.Cells.Item(20, 5).Value = 1
$wsLine.Cells.Item(20, 6).Value = 1

$wsLine.Cells.Item(21, 1).Value = 20
$wsLine.Cells.Item(21, 2).Value = 4
$wsLine.Cells.Item(21, 3).Value = 20
$wsLine.Cells.Item(21, 4).Value = 21
$wsLine.Cells.Item(21, 5).Value = 1
$wsLine.Cells.Item(21, 6).Value = 1

$wsLine.Range("H34").Select()

# ------------------------------------------------------------------
# 4) Nodal Load sheet: renumber the loaded nodes, switch the Fy load
#    value (as a shared formula chased down from D2) and extend the
#    table up to node 20.
# ------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Nodal Load")

$bVals = @(1,2,3,4,5,6,7,8,9,10,12,13,14,15,16,17,18,19,20)
for ($i = 0; $i -lt $bVals.Length; $i++) {
  $wsLoad.Cells.Item(2 + $i, 2).Value = $bVals[$i]
}

for ($r = 2; $r -le 20; $r++) {
  $wsLoad.Cells.Item($r, 1).Value = $r - 1
  $wsLoad.Cells.Item($r, 3).Value = 0
  $wsLoad.Cells.Item($r, 5).Value = 0
  $wsLoad.Cells.Item($r, 6).Value = 0
  $wsLoad.Cells.Item($r, 7).Value = 0
  $wsLoad.Cells.Item($r, 8).Value = 0
}

$wsLoad.Range("D2").Value = -20
$wsLoad.Range("D3").Formula = "=D2"
$wsLoad.Range("D4:D20").Formula = "=D3"

$wsLoad.Range("D3").Select()

# ------------------------------------------------------------------
# 5) Sections sheet: add the new "Moment Curvature Curve Sheet Name"
#    column, pointing each section at the "curve1" sheet. The
#    "curve1" shared string is interned here (before the curve1
#    sheet's own headers) so the shared-strings table keeps the same
#    ordering as the source workbook.
# ------------------------------------------------------------------
$wsSect = $wb.Worksheets.Item("Sections")

$wsSect.Range("H2").Value = "curve1"

# ------------------------------------------------------------------
# 6) New "curve1" sheet holding the moment-curvature pairs, inserted
#    after "Nodal Load" (the last existing sheet).
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCurve = $wb.Worksheets.Add($null, $lastSheet)
$wsCurve.Name = "curve1"

$wsCurve.Range("A1").Value = "Curvature"
$wsCurve.Range("B1").Value = "Moment"

$curveData = @(
  @(-0.39660499999999999, -4.1699999999999999),
  @(-0.32197399999999998, -4.7800000000000002),
  @(-0.228134, -5.1399999999999997),
  @(-0.20505499999999999, -5.25),
  @(-0.181115, -5.3399999999999999),
  @(-0.15625500000000001, -5.4100000000000001),
  @(-0.12990199999999999, -5.4100000000000001),
  @(-0.12650500000000001, -5.4000000000000004),
  @(-0.12559000000000001, -5.4100000000000001),
  @(-0.125198, -5.4000000000000004),
  @(-0.12474, -5.4100000000000001),
  @(-0.12389, -5.4000000000000004),
  @(-0.12304, -5.4000000000000004),
  @(-0.116037, -5.3799999999999999),
  @(-0.10882600000000001, -5.3499999999999996),
  @(-0.10147399999999999, -5.3099999999999996),
  @(-0.085951, -5.2000000000000002),
  @(-0.069375000000000006, -5.0300000000000002),
  @(-0.051763000000000003, -4.7300000000000004),
  @(-0.033103, -4.4100000000000001),
  @(-0.023764, -4.2000000000000002),
  @(-0.014063000000000001, -3.96),
  @(-0.0067889999999999999, -2.21),
  @(0, 0),
  @(0.006705, 9.1099999999999994),
  @(0.013016, 17.34),
  @(0.014540000000000001, 19.260000000000002),
  @(0.016641, 20.449999999999999),
  @(0.017405, 20.5),
  @(0.018169000000000001, 20.539999999999999),
  @(0.018534999999999999, 20.57),
  @(0.018932999999999998, 20.579999999999998),
  @(0.019299, 20.609999999999999),
  @(0.019696999999999999, 20.59),
  @(0.022852000000000001, 20.789999999999999),
  @(0.037246000000000001, 21.129999999999999),
  @(0.044734000000000003, 21.219999999999999),
  @(0.052347999999999999, 21.300000000000001),
  @(0.059759, 21.34),
  @(0.066999000000000003, 21.350000000000001),
  @(0.080966999999999997, 21.370000000000001),
  @(0.094085000000000002, 21.329999999999998),
  @(0.140766, 21.120000000000001),
  @(0.17944099999999999, 20.800000000000001),
  @(0.25011699999999998, 18.870000000000001)
)

for ($i = 0; $i -lt $curveData.Length; $i++) {
  $row = 2 + $i
  $wsCurve.Cells.Item($row, 1).Value = $curveData[$i][0]
  $wsCurve.Cells.Item($row, 2).Value = $curveData[$i][1]
}

$wsCurve.Range("K19").Select()

# ------------------------------------------------------------------
# 7) Finish the Sections sheet header + selection now that the
#    earlier shared strings have been interned in the right order.
# ------------------------------------------------------------------
$wsSect.Range("H1").Value = "Moment Curvature Curve Sheet Name"

$wsSect.Range("H9").Select()
